$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all data rows up by one: removes the old first data row (2007/39400)
# and moves every subsequent row up, which also drops the old last row (19)
# and updates the dimension from A1:E19 to A1:E18.
$ws.Rows("2").Delete()

# Recalculated y_1_forecast (column E) values for each remaining data row.
$eValues = @{
    2  = -4.700509864312973
    3  = -0.01655958389530365
    4  = 3.579142225970444
    5  = -0.289184878867832
    6  = 5.963492031746176
    7  = 7.523777575896196
    8  = 2.532215190177589
    9  = 2.051185924063259
    10 = 0.4575538530338541
    11 = 2.600569166164624
    12 = 3.605726003451304
    13 = 3.490656491795074
    14 = -2.347097924577757
    15 = -0.1803381976702711
    16 = -1.152671696465724
    17 = -2.785556326028149
    18 = -2.452009576682213
}

foreach ($row in $eValues.Keys) {
    $ws.Cells.Item($row, 5).Value = $eValues[$row]
}
